$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Tipo shifts from D to E)
$ws.Columns.Item(4).Insert()

# Set the new header "MAE" in D1 and copy style from C1 (header style)
$ws.Range("D1").Value = "MAE"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Fill MAE values
$ws.Range("D2").Value = 0.8583964064824784
$ws.Range("D3").Value = 1.863313287418611
$ws.Range("D4").Value = 1.29915171380136
